$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '26.941.34'
Set-TextValue "E2" '  -1.56%  '

# Row 3
Set-TextValue "D3" '1.817.92'
Set-TextValue "E3" '  -0.84%  '

# Row 4
Set-TextValue "D4" '1.009'
Set-TextValue "E4" '  -0.13%  '

# Row 5
Set-TextValue "B5" 'USDC'
Set-TextValue "C5" 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue "D5" '1.008'
Set-TextValue "E5" '  -0.19%  '

# Row 6
Set-TextValue "B6" 'BNB'
Set-TextValue "C6" 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue "D6" '308.92'
Set-TextValue "E6" '  -1.78%  '

# Row 7
Set-TextValue "D7" '0.4662'
Set-TextValue "E7" '  -1.70%  '

# Row 8
Set-TextValue "D8" '0.3653'
Set-TextValue "E8" '  -1.15%  '

# Row 9
Set-TextValue "D9" '0.07223'
Set-TextValue "E9" '  -3.14%  '

# Row 10
Set-TextValue "E10" '  -3.19%  '

# Row 11
Set-TextValue "D11" '19.74'
Set-TextValue "E11" '  -3.53%  '

# Row 12
Set-TextValue "D12" '0.07554'
Set-TextValue "E12" '  +3.01%  '

# Row 13
Set-TextValue "B13" 'Polkadot'
Set-TextValue "C13" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D13" '5.320'
Set-TextValue "E13" '  -2.46%  '

# Row 14
Set-TextValue "B14" 'Litecoin'
Set-TextValue "C14" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D14" '91.73'
Set-TextValue "E14" '  -1.46%  '

# Row 15
Set-TextValue "B15" 'Chainlink'
Set-TextValue "C15" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D15" '6.479'
Set-TextValue "E15" '  -1.65%  '

# Row 16
Set-TextValue "B16" 'WrappedEther'
Set-TextValue "C16" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D16" '1.671.36'
Set-TextValue "E16" '  -12.17%  '

# Row 17
Set-TextValue "D17" '1.009'
Set-TextValue "E17" '  -0.05%  '

# Row 18
Set-TextValue "D18" '0.000008616'
Set-TextValue "E18" '  -2.31%  '

# Row 19
Set-TextValue "D19" '1.008'
Set-TextValue "E19" '  -0.29%  '

# Row 21
Set-TextValue "D21" '26.843.05'
Set-TextValue "E21" '  -2.92%  '

# Row 22
Set-TextValue "D22" '5.135'
Set-TextValue "E22" '  -3.34%  '

# Row 23
Set-TextValue "D23" '10.52'
Set-TextValue "E23" '  -1.54%  '

# Row 24
Set-TextValue "D24" '1.943.56'
Set-TextValue "E24" '  -7.99%  '

# Row 25
Set-TextValue "D25" '151.78'
Set-TextValue "E25" '  -0.12%  '

# Row 26
Set-TextValue "E26" '  -2.32%  '

# Row 27
Set-TextValue "D27" '18.12'
Set-TextValue "E27" '  -2.77%  '

# Row 28
Set-TextValue "D28" '2.064'
Set-TextValue "E28" '  -3.66%  '

# Row 29
Set-TextValue "D29" '5.094'
Set-TextValue "E29" '  -2.97%  '

# Row 30
Set-TextValue "D30" '115.24'
Set-TextValue "E30" '  -1.85%  '

# Row 31
Set-TextValue "D31" '0.08869'
Set-TextValue "E31" '  -1.47%  '

# Row 32
Set-TextValue "D32" '2.962'
Set-TextValue "E32" '  +0.59%  '

# Row 33
Set-TextValue "D33" '4.414'
Set-TextValue "E33" '  -3.00%  '

# Row 34
Set-TextValue "E34" '  -4.35%  '

# Row 35
Set-TextValue "D35" '0.7155'
Set-TextValue "E35" '  -5.27%  '

# Row 36
Set-TextValue "D36" '1.077'
Set-TextValue "E36" '  -2.37%  '

# Row 37
Set-TextValue "D37" '0.05247'
Set-TextValue "E37" '  -2.03%  '

# Row 38
Set-TextValue "B38" 'VeChain'
Set-TextValue "C38" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D38" '0.01920'
Set-TextValue "E38" '  -1.70%  '

# Row 39
Set-TextValue "B39" 'RenderToken'
Set-TextValue "C39" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D39" '2.389'
Set-TextValue "E39" '  -0.40%  '

# Row 40
Set-TextValue "E40" '  -2.24%  '

# Row 41
Set-TextValue "D41" '7.130'
Set-TextValue "E41" '  -2.60%  '

# Row 42
Set-TextValue "E42" '  -3.43%  '

# Row 43
Set-TextValue "D43" '0.1621'
Set-TextValue "E43" '  -2.38%  '

# Row 44
Set-TextValue "D44" '8.141'
Set-TextValue "E44" '  -4.18%  '

# Row 45
Set-TextValue "D45" '0.4807'
Set-TextValue "E45" '  -2.22%  '

# Row 46
Set-TextValue "E46" '  -0.25%  '

# Row 47
Set-TextValue "D47" '10.09'
Set-TextValue "E47" '  -4.47%  '

# Row 48
Set-TextValue "D48" '102.86'
Set-TextValue "E48" '  -2.16%  '

# Row 49
Set-TextValue "D49" '0.06250'
Set-TextValue "E49" '  -0.79%  '

# Row 50
Set-TextValue "D50" '1.613'
Set-TextValue "E50" '  -3.77%  '

# Row 51
Set-TextValue "D51" '64.07'
Set-TextValue "E51" '  -2.61%  '
